$wb = $excel.ActiveWorkbook

# Update "10000" sheet (sheet1)
$ws = $wb.Worksheets.Item("10000")
$ws.Range("B2").Value = 3500
$ws.Range("C2").Value = 35
$ws.Range("B3").Value = 531
$ws.Range("C3").Value = 5.31
$ws.Range("B4").Value = 2110
$ws.Range("C4").Value = 21.1
$ws.Range("B5").Value = 3859
$ws.Range("C5").Value = 38.59
$ws.Range("C6").Value = 52.7

# Update "100000" sheet (sheet2)
$ws = $wb.Worksheets.Item("100000")
$ws.Range("B2").Value = 33976
$ws.Range("C2").Value = 33.976
$ws.Range("B3").Value = 6110
$ws.Range("C3").Value = 6.11
$ws.Range("B4").Value = 20817
$ws.Range("C4").Value = 20.817
$ws.Range("B5").Value = 39097
$ws.Range("C5").Value = 39.097
$ws.Range("C6").Value = 51.86

# Update "100000000" sheet (sheet3)
$ws = $wb.Worksheets.Item("100000000")
$ws.Range("B2").Value = 34006624
$ws.Range("C2").Value = 34.006624
$ws.Range("B3").Value = 6000263
$ws.Range("C3").Value = 6.000263
$ws.Range("B4").Value = 21003694
$ws.Range("C4").Value = 21.003694
$ws.Range("B5").Value = 38989419
$ws.Range("C5").Value = 38.989419
$ws.Range("C6").Value = 52

# Add new "1000000000" sheet (sheet4) after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "1000000000"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Header row, styled like the other sheets' header (bold, centered, top-aligned, thin border)
$newSheet.Range("A1").Value = "Final"
$newSheet.Range("B1").Value = "Ocorrências"
$newSheet.Range("C1").Value = "Probabilidade"
$headerRange = $newSheet.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Match page margins used by the other sheets (0.75in/1in/0.5in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$newSheet.Range("A2").Value = "Final A"
$newSheet.Range("B2").Value = 339996378
$newSheet.Range("C2").Value = 33.9996378

$newSheet.Range("A3").Value = "Final B"
$newSheet.Range("B3").Value = 60004236
$newSheet.Range("C3").Value = 6.0004236

$newSheet.Range("A4").Value = "Final C"
$newSheet.Range("B4").Value = 209975093
$newSheet.Range("C4").Value = 20.9975093

$newSheet.Range("A5").Value = "Final D"
$newSheet.Range("B5").Value = 390024293
$newSheet.Range("C5").Value = 39.0024293

$newSheet.Range("A6").Value = "Cutscene do Cavalo"
$newSheet.Range("B6").Value = ""
$newSheet.Range("C6").Value = 52

# Restore original active sheet/selection
$wb.Worksheets.Item("10000").Activate()
